# Populate the "defis" (challenges) sheet with its real content:
#  - column A: category name (already present, kept as-is)
#  - column B: long-form explanation / advice for that category
#  - column C (new): short, punchy "myth vs. reality" title for the category
# Row heights grow to fit the now much longer wrapped text in column B.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column C: short title/"myth" strings, right-hand side of the table.
$ws.Columns.Item(3).ColumnWidth = 30.21875

$ws.Range("A1").Value2 = 'Chaussures'
$ws.Range("B1").Value2 = 'Débarrasse-toi des chaussures inconfortables, qui te font mal aux pieds, trop abîmées pour être portées.<br>Ne garde qu''une paire de chaussures "pour les occasions".<br>Essaie de te limiter à 2 paires de chaussures par saison.'
$ws.Range("C1").Value2 = 'Je pourrais chausser un mille-pattes'
$ws.Rows.Item(1).RowHeight = 28.8

$ws.Range("A2").Value2 = 'Sous-vêtements'
$ws.Range("B2").Value2 = 'Jette tes sous-vêtements détendus, troués, inconfortables, trop petits ou trop grands…<br>Une astuce pour faire le tri : si tu meurs de honte à l''idée qu''on puisse te voir dans tel ou tel sous-vêtement, ne le garde pas !'
$ws.Range("C2").Value2 = 'Imagine que tu te retrouves soudainement en sous-vêtements dans la rue #monpirecauchemar'
$ws.Rows.Item(2).RowHeight = 28.8

$ws.Range("A3").Value2 = 'Vêtements du quotidien'
$ws.Range("B3").Value2 = 'Débarrasse-toi <br>- des vêtements inconfortables, dans lesquels tu ne te sens pas bien,<br>- des vêtements trop abîmés pour que tu oses les porter,<br>- des vêtements trop petits et trop grands,<br>- des vêtements qui ne vont pas avec le reste de ta garde-robe (couleurs, coupes, style...),<br>- des vêtements que tu ne portes pas régulièrement (excepté les vêtements pour les occasions spéciales).<br><br><a href="https://www.theflonicles.be/2019/04/pourquoi-comment-garde-robe-capsule.html">Comment créer une garde-robe minimaliste pour s''habiller en 5 minutes les yeux fermés ?</a>'
$ws.Range("C3").Value2 = 'Plus de vêtements, moins de choix !'
$ws.Rows.Item(3).RowHeight = 72

$ws.Range("A4").Value2 = 'Vêtements de sport'
$ws.Range("B4").Value2 = 'Garde le nombre de vêtements nécessaires en fonction du sport que tu pratiques effectivement et en tenant compte de ton rythme de lessives. Posséder 3 paires de baskets de running ne te motivera pas à courir plus souvent !'
$ws.Range("C4").Value2 = 'Plus j''en ai, plus je serai motivé.e, non ?'
$ws.Rows.Item(4).RowHeight = 28.8

$ws.Range("A5").Value2 = 'Vêtements maison et nuit'
$ws.Range("B5").Value2 = 'Il s''agit de tous les vêtements abîmés ou un peu vieux, sans forme, que tu gardes pour rester à la maison, pour dormir ou pour les travaux salissants, ainsi que les vêtements achetés dans ce but. Débarrasse-toi du superflu pour ne garder que ce dont tu as vraiment besoin : <br>- pas la peine de garder 10 t-shirts pour les travaux manuels si tu n''as pas de projet concret dans les mois à venir,<br>- tes vêtements vont continuer de s''abîmer et tu ne risques pas de manquer de t-shirts défraîchis et déformés dans le futur,<br>- limite-toi au nombre de vêtements nécessaires pour tenir entre deux lessives.'
$ws.Range("C5").Value2 = 'Plus c''est détendu, mieux c''est'
$ws.Rows.Item(5).RowHeight = 72

$ws.Range("A6").Value2 = 'Manteaux et vestes'
$ws.Range("B6").Value2 = 'Garde un manteau ou une veste pour chaque saison, qui s''accorde avec toute ta garde-robe. Mise sur un classique intemporel.'
$ws.Range("C6").Value2 = 'Un manteau pour chaque jour de la semaine'
$ws.Rows.Item(6).RowHeight = 28.8

$ws.Range("A7").Value2 = 'Bijoux et accessoires'
$ws.Range("B7").Value2 = 'Conserve 2 ou 3 bijoux de chaque type maximum, qui s''accordent avec la majeure partie de ta garde-robe.<br>Ne garde que les accessoires nécessaires, en 1 ou 2 exemplaires chacun (ceintures, sacs, cravates/nœuds, bretelles…).'
$ws.Range("C7").Value2 = 'Keep it simple!'
$ws.Rows.Item(7).RowHeight = 28.8

$ws.Range("A8").Value2 = 'Vêtements pour les occasions'
$ws.Range("B8").Value2 = 'Garde une tenue pour les "grandes occasions", dans laquelle tu te sens vraiment bien. Si tu as vraiment "besoin" d''une tenue différente pour l''un ou l''autre événement, tu peux toujours l''emprunter ou la louer. Du reste, porter la même tenue d''un mariage à l''autre n''est pas une catastrophe !<br><br>Fais également le tri dans les tenues que tu gardes pour les vacances pour ne garder que le nombre nécessaire. Mieux : pioche dans ta garde-robe du quotidien lorsque tu fais ta valise et évite ainsi de laisser des vêtements dormir dans une armoire pendant un an.'
$ws.Range("C8").Value2 = 'Ma tenue signature'
$ws.Rows.Item(8).RowHeight = 72

$ws.Range("A9").Value2 = 'Vêtements sentimentaux'
$ws.Range("B9").Value2 = 'Les t-shirts de concert, souvenirs de la vie estudiantine ou d''un voyage… S''ils restent enfermés dans une boîte et que tu ne les portes pas, te servent-ils à quelque chose ?'
$ws.Range("C9").Value2 = 'Les souvenirs, c''est dans la tête'
$ws.Rows.Item(9).RowHeight = 28.8

$ws.Range("A10").Value2 = 'Linge de lit'
$ws.Range("B10").Value2 = 'Trie ton linge de lit pour te débarrasser de ce qui est en très mauvais état, des linges dépareillés, ceux qui grattent ou qui ont mal vieilli, pour ne garder que ceux qui sont confortables et que tu utilises. Conserve le nombre nécessaire pour pouvoir tenir entre deux lessives.'
$ws.Rows.Item(10).RowHeight = 43.2

# Restore the wrap-text style on column B (explanatory text) for every row —
# it already carried style index 1 on most rows; make sure it is applied
# uniformly now that every row has long text in B.
$ws.Range("B1:B10").WrapText = $true

# Recreate the saved view state: scrolled so row 2 is at the top, with the
# last-edited cell (C10) selected.
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C10").Select() | Out-Null
